$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text content parses as a plain decimal number (single dot).
# Force Text format first so Excel keeps them as literal strings instead of
# silently converting them to numeric cells (matches the source data, which
# stores every Price cell as inline text, even ones that look numeric).
$textCells = @("D5","D6","D10","D11","D14","D19","D20","D23","D25","D29","D33","D35","D38","D45","D47")
foreach ($cellref in $textCells) {
    $ws.Range($cellref).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.153.36"
$ws.Range("E2").Value = "  -2.22%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.574.51"
$ws.Range("E3").Value = "  -1.77%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.49%  "

# Row 5 - BNB
$ws.Range("D5").Value = "208.82"
$ws.Range("E5").Value = "  -1.47%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.498"
$ws.Range("E6").Value = "  -3.10%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.45%  "

# Row 8 - Dogecoin
$ws.Range("E8").Value = "  -1.63%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -1.27%  "

# Row 10 - Solana
$ws.Range("D10").Value = "19.53"
$ws.Range("E10").Value = "  -0.69%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.0843"
$ws.Range("E11").Value = "  -0.58%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.796.04"

# Row 13 - now WrappedEther (was Polkadot)
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.584.09"
$ws.Range("E13").Value = "  -0.76%  "

# Row 14 - now Polkadot (was WrappedEther)
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "4.06"
$ws.Range("E14").Value = "  -0.27%  "

# Row 16 - Litecoin
$ws.Range("E16").Value = "  -1.08%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.146.83"
$ws.Range("E17").Value = "  -2.18%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  -2.24%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "7.26"
$ws.Range("E19").Value = "  +1.80%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "207.63"
$ws.Range("E20").Value = "  -0.94%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.41%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -0.99%  "

# Row 23 - Toncoin
$ws.Range("D23").Value = "2.17"
$ws.Range("E23").Value = "  -2.55%  "

# Row 24 - Avalanche
$ws.Range("E24").Value = "  -2.62%  "

# Row 25 - Monero
$ws.Range("D25").Value = "143.79"
$ws.Range("E25").Value = "  +0.04%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  -1.58%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  -1.62%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "15.19"
$ws.Range("E29").Value = "  -1.16%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -0.40%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -1.52%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -2.08%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "2.99"
$ws.Range("E33").Value = "  +0.78%  "

# Row 34 - Maker
$ws.Range("D34").Value = "1.279.38"
$ws.Range("E34").Value = "  -0.71%  "

# Row 35 - ImmutableX
$ws.Range("D35").Value = "0.613"
$ws.Range("E35").Value = "  +3.58%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -1.47%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  -0.90%  "

# Row 38 - WEMIXToken
$ws.Range("D38").Value = "1.12"
$ws.Range("E38").Value = "  -10.10%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -2.40%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  -2.27%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  -0.43%  "

# Row 43 - MXToken
$ws.Range("E43").Value = "  -2.88%  "

# Row 44 - TrustWalletToken
$ws.Range("E44").Value = "  -2.20%  "

# Row 45 - Aave
$ws.Range("D45").Value = "62.35"
$ws.Range("E45").Value = "  -0.77%  "

# Row 46 - RocketPoolETH
$ws.Range("D46").Value = "1.709.27"
$ws.Range("E46").Value = "  -1.76%  "

# Row 47 - Quant
$ws.Range("D47").Value = "88.74"
$ws.Range("E47").Value = "  -1.95%  "

# Row 48 - RenderToken
$ws.Range("E48").Value = "  -3.03%  "

# Row 49 - BabyDogeCoin
$ws.Range("E49").Value = "  +1.33%  "

# Row 50 - Algorand
$ws.Range("E50").Value = "  -1.61%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  -1.67%  "
